# Hortaliza / Terminal La Palmera de La Serena - Pepino dulce
# Weekly update: insert a new week of price data (week of 2022-03-22)
# at the top of the time series (rows 279-281), pushing all the
# subsequent historical rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 279, shifting the existing
# rows 279:331 down to 282:334.
$ws.Range("A279:A281").EntireRow.Insert()

# New weekly records to populate the freshly inserted rows.
# Columns: A Mercado ID, B Mercado, C Región, D Fecha (serial), E Codreg,
#          F Categoría ID, G Categoría, H Variedad, I Calidad, J Volumen,
#          K Precio mínimo, L Precio máximo, M Precio promedio ponderado,
#          N Unidad de comercialización, O Origen, P Precio $/Kg,
#          Q Kg o Unidades, R Clasificación
$newRows = @(
  @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44642, 4, 100112043, "Pepino dulce", "Cultivar IV Región", "Primera", 440, 9500, 10000, 9750, "`$/bandeja 18 kilos", "Provincia de Limarí", 542, 18, "Hortaliza"),
  @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44642, 4, 100112043, "Pepino dulce", "Cultivar IV Región", "Segunda", 280, 7500, 8000, 7750, "`$/bandeja 18 kilos", "Provincia de Limarí", 431, 18, "Hortaliza"),
  @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44642, 4, 100112043, "Pepino dulce", "Cultivar IV Región", "Tercera", 200, 5500, 6000, 5750, "`$/bandeja 18 kilos", "Provincia de Limarí", 319, 18, "Hortaliza")
)

$startRow = 279
for ($i = 0; $i -lt $newRows.Length; $i++) {
  $r = $startRow + $i
  $vals = $newRows[$i]
  for ($c = 1; $c -le $vals.Length; $c++) {
    $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
  }
}
